# edit.ps1
#
# Applies the "Create LoggerSingleton and logging_lock" commit to plan.xlsx:
#   - refreshes the "Scheduled Time" timestamp (column F) for the existing
#     Instagram-mislabelled rows 21-28 from 2024-02-14 15:38 to the newer run
#     at 2024-02-15 00:07 (serial 45337.004861111112)
#   - appends 10 new scheduled Instagram posts (rows 29-38, Post ID 28-37)
#     using the same timestamp, content "TestFI 19".."TestFI 28" and image
#     paths "img22".."img31"
#   - updates the active selection to F42 and scrolls the sheet so row 14 is
#     at the top, matching where the author left the cursor after appending
#     the rows
#
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Scheduled Time for existing rows 21-28 to the new timestamp
$ws.Cells.Item(21,6).Value = 45337.004861111112
$ws.Cells.Item(22,6).Value = 45337.004861111112
$ws.Cells.Item(23,6).Value = 45337.004861111112
$ws.Cells.Item(24,6).Value = 45337.004861111112
$ws.Cells.Item(25,6).Value = 45337.004861111112
$ws.Cells.Item(26,6).Value = 45337.004861111112
$ws.Cells.Item(27,6).Value = 45337.004861111112
$ws.Cells.Item(28,6).Value = 45337.004861111112

# Append 10 new Instagram rows (28-37) with TestFI 19..28 / img22..31
$ws.Cells.Item(29,1).Value = 28
$ws.Cells.Item(29,2).Value = "Instagram"
$ws.Cells.Item(29,3).Value = "TestFI 19"
$ws.Cells.Item(29,4).Value = "img22"
$ws.Cells.Item(29,5).Value = "#new #tech #insta"
$ws.Cells.Item(29,6).Value = 45337.004861111112
$ws.Cells.Item(29,6).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Cells.Item(29,7).Value = "Scheduled"

$ws.Cells.Item(30,1).Value = 29
$ws.Cells.Item(30,2).Value = "Instagram"
$ws.Cells.Item(30,3).Value = "TestFI 20"
$ws.Cells.Item(30,4).Value = "img23"
$ws.Cells.Item(30,5).Value = "#new #tech #insta"
$ws.Cells.Item(30,6).Value = 45337.004861111112
$ws.Cells.Item(30,6).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Cells.Item(30,7).Value = "Scheduled"

$ws.Cells.Item(31,1).Value = 30
$ws.Cells.Item(31,2).Value = "Instagram"
$ws.Cells.Item(31,3).Value = "TestFI 21"
$ws.Cells.Item(31,4).Value = "img24"
$ws.Cells.Item(31,5).Value = "#new #tech #insta"
$ws.Cells.Item(31,6).Value = 45337.004861111112
$ws.Cells.Item(31,6).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Cells.Item(31,7).Value = "Scheduled"

$ws.Cells.Item(32,1).Value = 31
$ws.Cells.Item(32,2).Value = "Instagram"
$ws.Cells.Item(32,3).Value = "TestFI 22"
$ws.Cells.Item(32,4).Value = "img25"
$ws.Cells.Item(32,5).Value = "#new #tech #insta"
$ws.Cells.Item(32,6).Value = 45337.004861111112
$ws.Cells.Item(32,6).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Cells.Item(32,7).Value = "Scheduled"

$ws.Cells.Item(33,1).Value = 32
$ws.Cells.Item(33,2).Value = "Instagram"
$ws.Cells.Item(33,3).Value = "TestFI 23"
$ws.Cells.Item(33,4).Value = "img26"
$ws.Cells.Item(33,5).Value = "#new #tech #insta"
$ws.Cells.Item(33,6).Value = 45337.004861111112
$ws.Cells.Item(33,6).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Cells.Item(33,7).Value = "Scheduled"

$ws.Cells.Item(34,1).Value = 33
$ws.Cells.Item(34,2).Value = "Instagram"
$ws.Cells.Item(34,3).Value = "TestFI 24"
$ws.Cells.Item(34,4).Value = "img27"
$ws.Cells.Item(34,5).Value = "#new #tech #insta"
$ws.Cells.Item(34,6).Value = 45337.004861111112
$ws.Cells.Item(34,6).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Cells.Item(34,7).Value = "Scheduled"

$ws.Cells.Item(35,1).Value = 34
$ws.Cells.Item(35,2).Value = "Instagram"
$ws.Cells.Item(35,3).Value = "TestFI 25"
$ws.Cells.Item(35,4).Value = "img28"
$ws.Cells.Item(35,5).Value = "#new #tech #insta"
$ws.Cells.Item(35,6).Value = 45337.004861111112
$ws.Cells.Item(35,6).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Cells.Item(35,7).Value = "Scheduled"

$ws.Cells.Item(36,1).Value = 35
$ws.Cells.Item(36,2).Value = "Instagram"
$ws.Cells.Item(36,3).Value = "TestFI 26"
$ws.Cells.Item(36,4).Value = "img29"
$ws.Cells.Item(36,5).Value = "#new #tech #insta"
$ws.Cells.Item(36,6).Value = 45337.004861111112
$ws.Cells.Item(36,6).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Cells.Item(36,7).Value = "Scheduled"

$ws.Cells.Item(37,1).Value = 36
$ws.Cells.Item(37,2).Value = "Instagram"
$ws.Cells.Item(37,3).Value = "TestFI 27"
$ws.Cells.Item(37,4).Value = "img30"
$ws.Cells.Item(37,5).Value = "#new #tech #insta"
$ws.Cells.Item(37,6).Value = 45337.004861111112
$ws.Cells.Item(37,6).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Cells.Item(37,7).Value = "Scheduled"

$ws.Cells.Item(38,1).Value = 37
$ws.Cells.Item(38,2).Value = "Instagram"
$ws.Cells.Item(38,3).Value = "TestFI 28"
$ws.Cells.Item(38,4).Value = "img31"
$ws.Cells.Item(38,5).Value = "#new #tech #insta"
$ws.Cells.Item(38,6).Value = 45337.004861111112
$ws.Cells.Item(38,6).NumberFormat = "yyyy\-mm\-dd\ hh:mm:ss"
$ws.Cells.Item(38,7).Value = "Scheduled"

# Update the view: scroll so row 14 is the top-left visible row, and select F42
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 14
$win.ScrollColumn = 1
$ws.Range("F42").Select()
